# Scheduled runner: refresh market-board price/profit figures across all crafting-job sheets.
# For each touched leve row, currentAveragePrice / NQ / HQ prices and the derived Leve profit
# columns (H-N) are recomputed from freshly-fetched market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 425.54544
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 425.54544
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1276.63632
$ws.Range("N17").Value = -1612.63632

$ws.Range("H86").Value = 3777.7778
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 4750
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4750
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -6996

$ws.Range("H88").Value = 1832
$ws.Range("I88").Value = 969.6667
$ws.Range("J88").Value = 2694.3333
$ws.Range("K88").Value = 969.6667
$ws.Range("L88").Value = 2694.3333
$ws.Range("M88").Value = -563.6667
$ws.Range("N88").Value = -3506.3333

$ws.Range("H89").Value = 3777.7778
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 4750
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 23750
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -34982

$ws.Range("H91").Value = 1832
$ws.Range("I91").Value = 969.6667
$ws.Range("J91").Value = 2694.3333
$ws.Range("K91").Value = 969.6667
$ws.Range("L91").Value = 2694.3333
$ws.Range("M91").Value = 434.3333
$ws.Range("N91").Value = -5502.3333

$ws.Range("H98").Value = 10323.467
$ws.Range("I98").Value = 5961.3335
$ws.Range("J98").Value = 16866.666
$ws.Range("K98").Value = 5961.3335
$ws.Range("L98").Value = 16866.666
$ws.Range("M98").Value = -4463.3335
$ws.Range("N98").Value = -19862.666

$ws.Range("H116").Value = 16668886
$ws.Range("I116").Value = 200000000
$ws.Range("J116").Value = 2421.4546
$ws.Range("K116").Value = 200000000
$ws.Range("L116").Value = 2421.4546
$ws.Range("M116").Value = -199996558
$ws.Range("N116").Value = -9305.454600000001

$ws.Range("H122").Value = 10323.467
$ws.Range("I122").Value = 5961.3335
$ws.Range("J122").Value = 16866.666
$ws.Range("K122").Value = 17884.0005
$ws.Range("L122").Value = 50599.99800000001
$ws.Range("M122").Value = -15434.0005
$ws.Range("N122").Value = -55499.99800000001

$ws.Range("H129").Value = 1042.6
$ws.Range("I129").Value = 393.66666
$ws.Range("J129").Value = 1136.4216
$ws.Range("K129").Value = 1180.99998
$ws.Range("L129").Value = 3409.2648
$ws.Range("M129").Value = 3819.00002
$ws.Range("N129").Value = -13409.2648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4800
$ws.Range("I16").Value = 4800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4800
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4513

$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -1123

$ws.Range("H122").Value = 8812
$ws.Range("I122").Value = 8812
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 26436
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -23986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 29967.938
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 29967.938
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 29967.938
$ws.Range("N62").Value = -31339.938

$ws.Range("H65").Value = 29967.938
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 29967.938
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 89903.814
$ws.Range("N65").Value = -96767.814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4160.722
$ws.Range("I31").Value = 1631.1
$ws.Range("J31").Value = 5133.654
$ws.Range("K31").Value = 1631.1
$ws.Range("L31").Value = 5133.654
$ws.Range("M31").Value = -1336.1
$ws.Range("N31").Value = -5723.654

$ws.Range("H34").Value = 4160.722
$ws.Range("I34").Value = 1631.1
$ws.Range("J34").Value = 5133.654
$ws.Range("K34").Value = 1631.1
$ws.Range("L34").Value = 5133.654
$ws.Range("M34").Value = -1429.1
$ws.Range("N34").Value = -5537.654

$ws.Range("H122").Value = 1556
$ws.Range("I122").Value = 1556
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4668
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2218

$ws.Range("H132").Value = 2774.7827
$ws.Range("I132").Value = 2101.8333
$ws.Range("J132").Value = 3508.9092
$ws.Range("K132").Value = 6305.499899999999
$ws.Range("L132").Value = 10526.7276
$ws.Range("M132").Value = -3775.499899999999
$ws.Range("N132").Value = -15586.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 7550
$ws.Range("I17").Value = 4650
$ws.Range("J17").Value = 9000
$ws.Range("K17").Value = 13950
$ws.Range("L17").Value = 27000
$ws.Range("M17").Value = -13781
$ws.Range("N17").Value = -27338

$ws.Range("H121").Value = 1227.8
$ws.Range("I121").Value = 878.3333
$ws.Range("J121").Value = 1377.5714
$ws.Range("K121").Value = 2634.9999
$ws.Range("L121").Value = 4132.7142
$ws.Range("M121").Value = -1324.9999
$ws.Range("N121").Value = -6752.7142

$ws.Range("H131").Value = 20453.086
$ws.Range("I131").Value = 756
$ws.Range("J131").Value = 22797.977
$ws.Range("K131").Value = 2268
$ws.Range("L131").Value = 68393.931
$ws.Range("M131").Value = 2772
$ws.Range("N131").Value = -78473.931

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4527.8237
$ws.Range("I102").Value = 4472.8335
$ws.Range("J102").Value = 4659.8
$ws.Range("K102").Value = 4472.8335
$ws.Range("L102").Value = 4659.8
$ws.Range("M102").Value = -2850.8335
$ws.Range("N102").Value = -7903.8

$ws.Range("H131").Value = 49163
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 49163
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 49163
$ws.Range("N131").Value = -59243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4913.1816
$ws.Range("I7").Value = 4005.7144
$ws.Range("J7").Value = 6501.25
$ws.Range("K7").Value = 4005.7144
$ws.Range("L7").Value = 6501.25
$ws.Range("M7").Value = -3893.7144
$ws.Range("N7").Value = -6725.25

$ws.Range("H40").Value = 2977.8572
$ws.Range("I40").Value = 2974
$ws.Range("J40").Value = 3001
$ws.Range("K40").Value = 2974
$ws.Range("L40").Value = 3001
$ws.Range("M40").Value = -2838
$ws.Range("N40").Value = -3273

$ws.Range("H126").Value = 4913.1816
$ws.Range("I126").Value = 4005.7144
$ws.Range("J126").Value = 6501.25
$ws.Range("K126").Value = 12017.1432
$ws.Range("L126").Value = 19503.75
$ws.Range("M126").Value = -9547.143199999999
$ws.Range("N126").Value = -24443.75

$ws.Range("H132").Value = 5205.4443
$ws.Range("I132").Value = 6320
$ws.Range("J132").Value = 3812.25
$ws.Range("K132").Value = 18960
$ws.Range("L132").Value = 11436.75
$ws.Range("M132").Value = -16430
$ws.Range("N132").Value = -16496.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 39900
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 39900
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 39900
$ws.Range("N54").Value = -40940

$ws.Range("H123").Value = 24428.475
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24428.475
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24428.475
$ws.Range("N123").Value = -34228.475

$ws.Range("H124").Value = 200286
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 200286
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 200286
$ws.Range("N124").Value = -210106

$ws.Range("H126").Value = 7427.0415
$ws.Range("I126").Value = 9114.941000000001
$ws.Range("J126").Value = 3327.8572
$ws.Range("K126").Value = 27344.823
$ws.Range("L126").Value = 9983.571599999999
$ws.Range("M126").Value = -24874.823
$ws.Range("N126").Value = -14923.5716
